# Replace "NaN" text values with "NA" in the worksheet.
# The three cells B20, E21, B22 currently hold the string "NaN";
# the commit replaces them with the string "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B20").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("B22").Value = "NA"

# Reflect the resulting scroll position / selection left by editing B20.
$ws.Application.Goto($ws.Range("A7"))
$ws.Range("B20").Select()
